# Update "想去人数" (want-to-go count) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 124
$ws1.Range("F4").Value = 161
$ws1.Range("F5").Value = 3173
$ws1.Range("F7").Value = 9

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 124
$ws4.Range("F4").Value = 161
$ws4.Range("F5").Value = 3173
$ws4.Range("F9").Value = 9
